$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 85; existing rows 85:178 shift down to 86:179.
$ws.Rows("85:85").Insert()

# Populate the newly inserted row 85 with the new data record.
$ws.Range("A85").Value2 = 9
$ws.Range("B85").Value2 = "Vega Central Mapocho de Santiago"
$ws.Range("C85").Value2 = "Metropolitana"
$ws.Range("D85").Value2 = 44494
$ws.Range("E85").Value2 = 13
$ws.Range("F85").Value2 = 100112021
$ws.Range("G85").Value2 = "Ají"
$ws.Range("H85").Value2 = "Inferno"
$ws.Range("I85").Value2 = "Primera"
$ws.Range("J85").Value2 = 40
$ws.Range("K85").Value2 = 43000
$ws.Range("L85").Value2 = 45000
$ws.Range("M85").Value2 = 43800
$ws.Range("N85").Value2 = "`$/caja 12 kilos"
$ws.Range("O85").Value2 = "Región de Arica y Parinacota"
$ws.Range("P85").Value2 = 3650
$ws.Range("Q85").Value2 = 12
$ws.Range("R85").Value2 = "Hortaliza"
